$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New 2D training schedule data (rows 2-6, columns A-I); column J keeps the
# existing shared string "train_dim2_1" already present on each row.
$data = @(
    @(1, 2, 4, 7, 5, 5, 1, 12, 5),
    @(2, 0, 4, 4, 6, 4, 2, 23, 5),
    @(3, 4, 0, 5, 5, 1, 5, 56, 5),
    @(4, 3, 3, 6, 6, 3, 3, 34, 5),
    @(5, 1, 1, 3, 5, 2, 4, 45, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

$wb.Save()
